$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date the record was last changed, stored
# as an Excel serial date number. Every populated data row (2 through 111)
# currently shows 45190 (2023-09-21) and needs to be updated to 45192
# (2023-09-23).
for ($row = 2; $row -le 111; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
